$d = $word.ActiveDocument

# 1) Move the "_GoBack" bookmark from the end of the Sites table's Url row
#    (which is being deleted) to the start of the Pages table's
#    FoundDateTime value cell ("<<10-10-2015 15:23>>").
$pagesTable = $d.Tables.Item(4)
$foundDateCell = $pagesTable.Rows.Item(6).Cells.Item(4)
$bmRange = $d.Range($foundDateCell.Range.Start, $foundDateCell.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 2) Remove the "Url" row from the Sites table (last table in the document).
$sitesTable = $d.Tables.Item(5)
$sitesTable.Rows.Item(5).Delete()
